$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format before assigning values so that
# numeric-looking strings (e.g. "0.9964") are not auto-converted to numbers,
# matching the original inline-string cell type.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "22.357.09"
$ws.Range("E2").Value = "  +8.68%  "
$ws.Range("D3").Value = "1.599.96"
$ws.Range("E3").Value = "  +8.08%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "0.9964"
$ws.Range("E5").Value = "  +2.00%  "
$ws.Range("D6").Value = "289.33"
$ws.Range("E6").Value = "  +3.40%  "
$ws.Range("D7").Value = "0.3700"
$ws.Range("E7").Value = "  +0.96%  "
$ws.Range("D8").Value = "0.3403"
$ws.Range("E8").Value = "  +10.58%  "
$ws.Range("D9").Value = "42.66"
$ws.Range("E9").Value = "  +7.76%  "
$ws.Range("D10").Value = "1.143"
$ws.Range("E10").Value = "  +7.21%  "
$ws.Range("D11").Value = "0.07066"
$ws.Range("E11").Value = "  +5.98%  "
$ws.Range("E12").Value = "  -0.10%  "
$ws.Range("D13").Value = "19.81"
$ws.Range("E13").Value = "  +8.84%  "
$ws.Range("D14").Value = "5.944"
$ws.Range("E14").Value = "  +7.74%  "
$ws.Range("D15").Value = "6.667"
$ws.Range("E15").Value = "  +7.32%  "
$ws.Range("D16").Value = "0.00001087"
$ws.Range("E16").Value = "  +5.18%  "
$ws.Range("D17").Value = "0.9958"
$ws.Range("E17").Value = "  +1.99%  "
$ws.Range("D18").Value = "1.595.62"
$ws.Range("E18").Value = "  +8.00%  "
$ws.Range("D19").Value = "0.06618"
$ws.Range("E19").Value = "  +11.35%  "
$ws.Range("D20").Value = "78.46"
$ws.Range("E20").Value = "  +12.09%  "
$ws.Range("E21").Value = "  +11.30%  "
$ws.Range("D22").Value = "6.044"
$ws.Range("E22").Value = "  +10.05%  "
$ws.Range("D23").Value = "11.82"
$ws.Range("E23").Value = "  +6.80%  "
$ws.Range("D24").Value = "22.367.98"
$ws.Range("E24").Value = "  +8.55%  "
$ws.Range("D25").Value = "2.390"
$ws.Range("E25").Value = "  +6.47%  "
$ws.Range("D26").Value = "2.509"
$ws.Range("E26").Value = "  +16.95%  "
$ws.Range("D27").Value = "151.07"
$ws.Range("E27").Value = "  +7.01%  "
$ws.Range("D28").Value = "19.65"
$ws.Range("E28").Value = "  +13.68%  "
$ws.Range("D29").Value = "1.774.32"
$ws.Range("E29").Value = "  +8.56%  "
$ws.Range("D30").Value = "120.87"
$ws.Range("E30").Value = "  +5.64%  "
$ws.Range("D31").Value = "4.164"
$ws.Range("E31").Value = "  +5.41%  "
$ws.Range("D32").Value = "6.046"
$ws.Range("E32").Value = "  +20.75%  "
$ws.Range("D33").Value = "0.9530"
$ws.Range("E33").Value = "  +16.25%  "
$ws.Range("D34").Value = "0.08266"
$ws.Range("E34").Value = "  +2.76%  "
$ws.Range("D35").Value = "1.613"
$ws.Range("E35").Value = "  +5.79%  "
$ws.Range("D36").Value = "5.356"
$ws.Range("E36").Value = "  +13.36%  "
$ws.Range("D37").Value = "8.681"
$ws.Range("E37").Value = "  +11.91%  "
$ws.Range("E38").Value = "  +12.62%  "
$ws.Range("D39").Value = "0.06172"
$ws.Range("E39").Value = "  +5.44%  "
$ws.Range("D40").Value = "1.241"
$ws.Range("E40").Value = "  +1.99%  "
$ws.Range("D41").Value = "0.02219"
$ws.Range("E41").Value = "  +8.38%  "
$ws.Range("D42").Value = "0.2035"
$ws.Range("E42").Value = "  +7.71%  "
$ws.Range("D43").Value = "0.9955"
$ws.Range("E43").Value = "  +2.02%  "
$ws.Range("D44").Value = "0.5941"
$ws.Range("E44").Value = "  +11.93%  "
$ws.Range("D45").Value = "13.16"
$ws.Range("E45").Value = "  +7.72%  "
$ws.Range("D46").Value = "3.674"
$ws.Range("E46").Value = "  +4.06%  "
$ws.Range("D47").Value = "0.5725"
$ws.Range("E47").Value = "  +10.08%  "
$ws.Range("D48").Value = "126.07"
$ws.Range("E48").Value = "  +5.51%  "
$ws.Range("D49").Value = "1.977"
$ws.Range("E49").Value = "  +9.69%  "
$ws.Range("D50").Value = "0.06832"
$ws.Range("E50").Value = "  +5.47%  "
$ws.Range("D51").Value = "73.91"
$ws.Range("E51").Value = "  +9.13%  "

# Restore default styling on column D (the Text number format above would
# otherwise persist as a new cell style); this keeps cells unstyled like the rest
# of the data cells in the sheet.
$ws.Range("D2:D51").Style = "Normal"
